# Apply the delta report restructuring to the active worksheet.
#
# Old layout (A1:J4): Invoice No | Vendor Name | Invoice Date | GSTIN | PAN |
#                      HSN Code | Taxable Value | Total Amount | Status | Reason
# New layout (A1:F6):  Invoice No | Date | Vendor | GSTIN | Amount | Validation Status
# plus 5 fresh data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the columns that no longer exist (G:J) entirely, content + formatting ---
$ws.Range("G1:J4").Clear()

# --- Header row (A1:F1 keep their existing header style; just retext/reorder them) ---
$ws.Range("A1").Value = "Invoice No"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Vendor"
$ws.Range("D1").Value = "GSTIN"
$ws.Range("E1").Value = "Amount"
$ws.Range("F1").Value = "Validation Status"

# --- Clear old body content below the header so stale cells don't linger ---
$ws.Range("A2:F6").ClearContents()

# --- Data rows ---
$data = @(
    @("INV-1", "2025-06-14", "Vendor 1", "29ABCDE1234F01Z5", 1180, "VALID"),
    @("INV-2", "2025-06-14", "Vendor 2", "29ABCDE1234F02Z5", 1430, "VALID"),
    @("INV-3", "2025-06-14", "Vendor 3", "29ABCDE1234F03Z5", 1680, "VALID"),
    @("INV-4", "2025-06-14", "Vendor 4", "29ABCDE1234F04Z5", 1930, "VALID"),
    @("INV-5", "2025-06-14", "Vendor 5", "29ABCDE1234F05Z5", 2180, "VALID")
)

# Force column B (Date) to be stored as plain text, not an auto-converted date serial.
# Must happen before the values are assigned.
$ws.Range("B2:B6").NumberFormat = "@"

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $rowIndex++
}

# The forced "@" text format leaves column B with a distinct style; restore it
# back to the plain/default look (matching the rest of the data cells) now
# that the text values are safely committed.
$ws.Range("B2:B6").Style = $ws.Range("A2").Style
